$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the AutoFilter to the original 8-row table range first.
$ws.Range("A1:C8").AutoFilter()

# Rewrite the data rows 2-8 and append 4 new rows (9-12) with the final,
# already-sorted content (sorted ascending by the Cost column).
$data = @(
    @("Tactic 0", 1, "Put 1 card from the pool to the table discard pile"),
    @("Bench replacement", 3, "Search your deck and pick one footballer. Then, put 1 footballer from your hand to the deck. Shuffle your deck."),
    @("Tactic 3", 4, "Remove 1 card from your discard pile from the game OR draw a card"),
    @("4-3-3", 4, "You may roll +1 die to your defence tests this turn"),
    @("5-3-2", 4, "You may roll +1 die to you shooting tests this turn"),
    @("3-5-2", 4, "You may roll +1 die to you assisting tests this turn"),
    @("Tactic 4", 5, "You may skip a difficulty 3 or lower test of an event card that has two or more tests."),
    @("Tactic 1", 6, "You may re-roll 1 die for each test you do this turn"),
    @("Tactic 2", 6, "You may add 1 die for each test you do this turn"),
    @("Tactic 5", 6, "Gain 1 gold per test passed on this turn (max: 4)"),
    @("Coach anger", 6, "If you fail a test, play this card. You may try to pass the test again")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# The new rows need the same fill/font formatting as the rest of the table.
$ws.Range("A8:C8").Copy()
$ws.Range("A9:C12").PasteSpecial(-4122)

# Record the sort (by Cost, ascending) that produced this row order, now
# covering the full A2:C12 range.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B1:B8"))
$ws.Sort.SetRange($ws.Range("A2:C12"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Register the (hidden) _FilterDatabase name that Excel creates for the
# worksheet's AutoFilter.
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet3!`$A`$1:`$C`$8")
$n.Visible = $false

# Restore the active selection to the new last cell in column C.
$ws.Range("C11").Select()
